# Weekly refresh of the Fruta/Hortaliza "Higo" daily-logic subset:
# each row's Fecha/Volumen/Precio fields (and, for a few rows, Origen)
# are updated to reflect the newer source pull. Rows 12-13 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Higo, Primera)
$ws.Range("D2").Value = 44302
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 2143

# Row 3 (Higo, Segunda)
$ws.Range("D3").Value = 44302
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 1714

# Row 4 (Higo, Primera)
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 1714

# Row 5 (Higo, Segunda)
$ws.Range("D5").Value = 44322
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1143

# Row 6 (Higo, Primera)
$ws.Range("D6").Value = 44299
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("R6").Value = "Provincia de Santiago"
$ws.Range("S6").Value = 2143

# Row 7 (Higo, Segunda)
$ws.Range("D7").Value = 44299
$ws.Range("M7").Value = 75
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("R7").Value = "Provincia de Santiago"
$ws.Range("S7").Value = 1714

# Row 8 (Higo, Primera)
$ws.Range("D8").Value = 44320
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 1714

# Row 9 (Higo, Segunda)
$ws.Range("D9").Value = 44320
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 1143

# Row 10 (Higo, Primera)
$ws.Range("D10").Value = 44300
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("S10").Value = 2143

# Row 11 (Higo, Segunda)
$ws.Range("D11").Value = 44300
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("S11").Value = 1714

# Row 14 (Higo, Primera)
$ws.Range("D14").Value = 44292
$ws.Range("M14").Value = 25
$ws.Range("N14").Value = 16000
$ws.Range("O14").Value = 16000
$ws.Range("P14").Value = 16000
$ws.Range("S14").Value = 2286

# Row 15 (Higo, Segunda)
$ws.Range("D15").Value = 44292
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("S15").Value = 2143
